$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-150 down to 52-151.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with its data.
$ws.Cells.Item(51, 1).Value  = 5
$ws.Cells.Item(51, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(51, 3).Value  = "Maule"
$ws.Cells.Item(51, 4).Value  = 44662
$ws.Cells.Item(51, 5).Value  = 7
$ws.Cells.Item(51, 6).Value  = 100112017
$ws.Cells.Item(51, 7).Value  = "Apio"
$ws.Cells.Item(51, 8).Value  = "Americana (o)"
$ws.Cells.Item(51, 9).Value  = "Primera"
$ws.Cells.Item(51, 10).Value = 500
$ws.Cells.Item(51, 11).Value = 7500
$ws.Cells.Item(51, 12).Value = 7500
$ws.Cells.Item(51, 13).Value = 7500
$ws.Cells.Item(51, 14).Value = "`$/docena de matas"
$ws.Cells.Item(51, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(51, 16).Value = 1250
$ws.Cells.Item(51, 17).Value = 6
$ws.Cells.Item(51, 18).Value = "Hortaliza"
